# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header row (bold, bordered, centered)
# by copying the format from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record for every data row (2-44) with the same values
# for all players: 102 wins, 60 losses, 0 ties.
$ws.Range("AD2:AD44").Value = 102
$ws.Range("AE2:AE44").Value = 60
$ws.Range("AF2:AF44").Value = 0

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-44"
